# Scraper refresh for línea 141: 30/12/2025 12:47:09
# - LP1912 sheet: 15 new rows appended (222-236)
# - 6203-6173 sheet: 3 new rows appended (32-34)
# - All three sheets get their "Última actualización" timestamp bumped
# - LP1912 and 6203-6173 also get their "Total filas" counter bumped

$wb = $excel.ActiveWorkbook

$tsNew = "Última actualización: 30/12/2025 12:47:09"

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = $tsNew
$ws1.Range("A3").Value = "Total filas: 235"

$lp1912Rows = @(
    @(222, "12:46:58", "12:55", "10_OLMOS",             9),
    @(223, "12:46:58", "12:56", "16_SANTA ANA",         10),
    @(224, "12:46:58", "13:02", "15_ABASTO",            16),
    @(225, "12:46:58", "13:04", "23_HERNANDEZ",         18),
    @(226, "12:46:58", "13:06", "16_P MOR-SANTA ANA",   20),
    @(227, "12:46:58", "13:08", "10_OLMOS",             22),
    @(228, "12:46:58", "13:19", "10_OLMOS",             33),
    @(229, "12:46:58", "13:26", "14_ABASTO",            40),
    @(230, "12:46:58", "13:34", "23_HERNANDEZ",         48),
    @(231, "12:46:58", "13:36", "15_ABASTO",            50),
    @(232, "12:46:58", "13:46", "17_ROMERO",             60),
    @(233, "12:46:58", "13:56", "16_P MOR-167 Y 521",   70),
    @(234, "12:46:58", "14:04", "17_ROMERO",             78),
    @(235, "12:46:58", "14:04", "23_HERNANDEZ",         78),
    @(236, "12:46:58", "14:21", "26_HERNANDEZ",         95)
)

foreach ($row in $lp1912Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = "LP1912"
    $ws1.Cells.Item($r, 7).Value = "30/12/2025"
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (timestamp only, no new rows)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = $tsNew

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = $tsNew
$ws3.Range("A3").Value = "Total filas: 33"

$ws3.Cells.Item(32, 2).Value = "30/12/2025"
$ws3.Cells.Item(32, 3).Value = "12:47:04"
$ws3.Cells.Item(32, 4).Value = "12:50"
$ws3.Cells.Item(32, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(32, 6).Value = 3
$ws3.Cells.Item(32, 7).Value = "L6203"

$ws3.Cells.Item(33, 2).Value = "30/12/2025"
$ws3.Cells.Item(33, 3).Value = "12:47:09"
$ws3.Cells.Item(33, 4).Value = "13:31"
$ws3.Cells.Item(33, 5).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(33, 6).Value = 44
$ws3.Cells.Item(33, 7).Value = "L6173"

$ws3.Cells.Item(34, 2).Value = "30/12/2025"
$ws3.Cells.Item(34, 3).Value = "12:47:09"
$ws3.Cells.Item(34, 4).Value = "14:09"
$ws3.Cells.Item(34, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(34, 6).Value = 82
$ws3.Cells.Item(34, 7).Value = "L6173"
